# Generate Report for Handback
# -----------------------------------------------------------------------
# This script mutates localization-status.xlsx to reflect that the
# 45578986-...md source file has been handed back (in both the zh-cn and
# de-de target languages): the status text changes, each language sheet
# grows a "Latest Target File" / "Latest Handback File" pair of links for
# rows 2 and 3, and the "Latest Handback DateTime" column is stamped with
# the handback timestamp.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$hyperlinkColor = 15570276   # BGR int for RGB FF6495ED - matches the workbook's existing HyperLink font color

function Style-AsHyperlink($range) {
    $range.Font.Underline = 2
    $range.Font.Color = $hyperlinkColor
}

# ------------------------------------------------------------------
# 1. Overview sheet: status text for the two handed-off files
# ------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# ------------------------------------------------------------------
# 2. zh-cn sheet
# ------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("B2").Value = $newStatus
$zhcn.Range("B3").Value = $newStatus

# Row 2: Latest Target File / Latest Handback File
$zhcn.Range("E2").Value = "45578986-a814-4c70-9611-db0c6edc09ea.md"
$zhcn.Hyperlinks.Add($zhcn.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/edcaa9fcc9ccea72a662d58cfca052126e60ddc5/e2e/45578986-a814-4c70-9611-db0c6edc09ea.md", "", "", "45578986-a814-4c70-9611-db0c6edc09ea.md")
Style-AsHyperlink $zhcn.Range("E2")

$zhcn.Range("F2").Value = "45578986-a814-4c70-9611-db0c6edc09ea.9f728e5ef7db08d5a80eb99c5c40fcec9261724c.zh-cn.xlf"
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7856ddbeba936fcf912c716f3a8752cfd6a0c4f7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/45578986-a814-4c70-9611-db0c6edc09ea.9f728e5ef7db08d5a80eb99c5c40fcec9261724c.zh-cn.xlf", "", "", "45578986-a814-4c70-9611-db0c6edc09ea.9f728e5ef7db08d5a80eb99c5c40fcec9261724c.zh-cn.xlf")
Style-AsHyperlink $zhcn.Range("F2")

# Row 2: Latest Handback DateTime
$zhcn.Range("G2").Value = "2016-03-08 06:54:00"

# Row 3: Latest Target File / Latest Handback File
$zhcn.Range("E3").Value = "45578986-a814-4c70-9611-db0c6edc09ea.md"
$zhcn.Hyperlinks.Add($zhcn.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/edcaa9fcc9ccea72a662d58cfca052126e60ddc5/e2e/45578986-a814-4c70-9611-db0c6edc09ea.md", "", "", "45578986-a814-4c70-9611-db0c6edc09ea.md")
Style-AsHyperlink $zhcn.Range("E3")

$zhcn.Range("F3").Value = "45578986-a814-4c70-9611-db0c6edc09ea.9f728e5ef7db08d5a80eb99c5c40fcec9261724c.zh-cn.xlf"
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7856ddbeba936fcf912c716f3a8752cfd6a0c4f7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/45578986-a814-4c70-9611-db0c6edc09ea.9f728e5ef7db08d5a80eb99c5c40fcec9261724c.zh-cn.xlf", "", "", "45578986-a814-4c70-9611-db0c6edc09ea.9f728e5ef7db08d5a80eb99c5c40fcec9261724c.zh-cn.xlf")
Style-AsHyperlink $zhcn.Range("F3")

# Row 3: Latest Handback DateTime
$zhcn.Range("G3").Value = "2016-03-08 06:54:00"

# ------------------------------------------------------------------
# 3. de-de sheet
# ------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("B2").Value = $newStatus
$dede.Range("B3").Value = $newStatus

# Row 2: Latest Target File / Latest Handback File
$dede.Range("E2").Value = "45578986-a814-4c70-9611-db0c6edc09ea.md"
$dede.Hyperlinks.Add($dede.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/edcaa9fcc9ccea72a662d58cfca052126e60ddc5/e2e/45578986-a814-4c70-9611-db0c6edc09ea.md", "", "", "45578986-a814-4c70-9611-db0c6edc09ea.md")
Style-AsHyperlink $dede.Range("E2")

$dede.Range("F2").Value = "45578986-a814-4c70-9611-db0c6edc09ea.9f728e5ef7db08d5a80eb99c5c40fcec9261724c.de-de.xlf"
$dede.Hyperlinks.Add($dede.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ee2e76536899a793202b1ad772d390e2a646df6b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/45578986-a814-4c70-9611-db0c6edc09ea.9f728e5ef7db08d5a80eb99c5c40fcec9261724c.de-de.xlf", "", "", "45578986-a814-4c70-9611-db0c6edc09ea.9f728e5ef7db08d5a80eb99c5c40fcec9261724c.de-de.xlf")
Style-AsHyperlink $dede.Range("F2")

# Row 2: Latest Handback DateTime
$dede.Range("G2").Value = "2016-03-08 06:54:16"

# Row 3: Latest Target File / Latest Handback File
$dede.Range("E3").Value = "45578986-a814-4c70-9611-db0c6edc09ea.md"
$dede.Hyperlinks.Add($dede.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/edcaa9fcc9ccea72a662d58cfca052126e60ddc5/e2e/45578986-a814-4c70-9611-db0c6edc09ea.md", "", "", "45578986-a814-4c70-9611-db0c6edc09ea.md")
Style-AsHyperlink $dede.Range("E3")

$dede.Range("F3").Value = "45578986-a814-4c70-9611-db0c6edc09ea.9f728e5ef7db08d5a80eb99c5c40fcec9261724c.de-de.xlf"
$dede.Hyperlinks.Add($dede.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ee2e76536899a793202b1ad772d390e2a646df6b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/45578986-a814-4c70-9611-db0c6edc09ea.9f728e5ef7db08d5a80eb99c5c40fcec9261724c.de-de.xlf", "", "", "45578986-a814-4c70-9611-db0c6edc09ea.9f728e5ef7db08d5a80eb99c5c40fcec9261724c.de-de.xlf")
Style-AsHyperlink $dede.Range("F3")

# Row 3: Latest Handback DateTime
$dede.Range("G3").Value = "2016-03-08 06:54:16"

Write-Host "Handback report generated."
